$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 136. This pushes the existing rows 136-160
# (and their formatting, e.g. the date style on column D) down to 137-161,
# matching the diff's net effect of a one-row-down shift for that block.
$ws.Rows("136").Insert()

# Populate the newly inserted row 136 with a new weekly price observation.
$ws.Cells.Item(136, 1).Value = 5
$ws.Cells.Item(136, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(136, 3).Value = "Maule"
$ws.Cells.Item(136, 4).Value = 44504
$ws.Cells.Item(136, 5).Value = 7
$ws.Cells.Item(136, 6).Value = 100112008
$ws.Cells.Item(136, 7).Value = "Coliflor"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 5000
$ws.Cells.Item(136, 11).Value = 600
$ws.Cells.Item(136, 12).Value = 600
$ws.Cells.Item(136, 13).Value = 600
$ws.Cells.Item(136, 14).Value = "`$/unidad"
$ws.Cells.Item(136, 15).Value = "Región del Maule"
$ws.Cells.Item(136, 16).Value = 600
$ws.Cells.Item(136, 17).Value = 1
$ws.Cells.Item(136, 18).Value = "Hortaliza"
